# Investor advisor template: random password for new users + default advisor view
#
# Users are now given a randomly generated password when their account is
# created, so the uploader no longer needs to supply one - drop the
# "Password" column header (and its value) from the template.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").ClearContents() | Out-Null

# Default advisor view: scroll the sheet so column C is the left-most visible
# column and land the active cell on H4.
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("H4").Select() | Out-Null
